# Update "想去人数" (F column) and "最低票价" (G column) figures on the
# "展览" and "全部类型" sheets to the latest scraped counts.
#
# Mapping of worksheet row -> column -> new value (both sheets share the
# exact same data and therefore the exact same updates):
#   Row  2: G -> 65
#   Row  3: F -> 42
#   Row  4: F -> 350
#   Row  6: F -> 410
#   Row  7: F -> 113
#   Row  8: F -> 118
#   Row 11: F -> 58
#   Row 12: F -> 114
#   Row 13: F -> 1098
#   Row 14: F -> 1453
#   Row 15: F -> 320
#   Row 20: F -> 54
#   Row 21: F -> 95
#   Row 22: F -> 250
#   Row 24: F -> 301
#   Row 25: F -> 1669
#   Row 29: F -> 626
#   Row 31: F -> 3921
#   Row 33: F -> 456
#   Row 35: F -> 993
#   Row 39: F -> 90

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Col = "G"; Value = 65 },
    @{ Row = 3;  Col = "F"; Value = 42 },
    @{ Row = 4;  Col = "F"; Value = 350 },
    @{ Row = 6;  Col = "F"; Value = 410 },
    @{ Row = 7;  Col = "F"; Value = 113 },
    @{ Row = 8;  Col = "F"; Value = 118 },
    @{ Row = 11; Col = "F"; Value = 58 },
    @{ Row = 12; Col = "F"; Value = 114 },
    @{ Row = 13; Col = "F"; Value = 1098 },
    @{ Row = 14; Col = "F"; Value = 1453 },
    @{ Row = 15; Col = "F"; Value = 320 },
    @{ Row = 20; Col = "F"; Value = 54 },
    @{ Row = 21; Col = "F"; Value = 95 },
    @{ Row = 22; Col = "F"; Value = 250 },
    @{ Row = 24; Col = "F"; Value = 301 },
    @{ Row = 25; Col = "F"; Value = 1669 },
    @{ Row = 29; Col = "F"; Value = 626 },
    @{ Row = 31; Col = "F"; Value = 3921 },
    @{ Row = 33; Col = "F"; Value = 456 },
    @{ Row = 35; Col = "F"; Value = 993 },
    @{ Row = 39; Col = "F"; Value = 90 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range($u.Col + $u.Row).Value = $u.Value
    }
}
